$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date in column C for rows 2-20 from 45190 to 45192
for ($row = 2; $row -le 20; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
